# "Add files via upload" - populate the new "Points" column (B) on the
# Results sheet for every player row (2-99). Most players scored 0 points;
# a handful of rows have non-zero totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Default every row in B2:B99 to 0 ...
$ws.Range("B2:B99").Value = 0

# ... then overwrite the few rows that actually have points.
$ws.Range("B2").Value = 3000
$ws.Range("B18").Value = 2500
$ws.Range("B28").Value = 1500

# Leave the selection where the author last left it.
$ws.Range("B32").Select()
